$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 35162
$ws.Range("E2").Value = 687076224789
$ws.Range("F2").Value = 9734232794
$ws.Range("G2").Value = 1.26128

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 1896.57
$ws.Range("E3").Value = 228059391070
$ws.Range("F3").Value = 10608584320
$ws.Range("G3").Value = 3.14884

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 1.001
$ws.Range("E4").Value = 85358927601
$ws.Range("F4").Value = 12373158796
$ws.Range("G4").Value = 0.0718

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 244.11
$ws.Range("E5").Value = 37550287562
$ws.Range("F5").Value = 473599097
$ws.Range("G5").Value = 4.02452

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "XRP"
$ws.Range("D6").Value = 0.64471
$ws.Range("E6").Value = 34501044249
$ws.Range("F6").Value = 1138004336
$ws.Range("G6").Value = 5.07296

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "USDC"
$ws.Range("D7").Value = 0.999987
$ws.Range("E7").Value = 24562938966
$ws.Range("F7").Value = 4012060204
$ws.Range("G7").Value = 0.02933

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "SOL"
$ws.Range("C8").Value = "Solana"
$ws.Range("D8").Value = 41.75
$ws.Range("E8").Value = 17470329514
$ws.Range("F8").Value = 931833435
$ws.Range("G8").Value = 0.81813

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "STETH"
$ws.Range("C9").Value = "Lido Staked Ether"
$ws.Range("D9").Value = 1888.27
$ws.Range("E9").Value = 16756138643
$ws.Range("F9").Value = 7340541
$ws.Range("G9").Value = 2.56979

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "ADA"
$ws.Range("C10").Value = "Cardano"
$ws.Range("D10").Value = 0.344995
$ws.Range("E10").Value = 12038944621
$ws.Range("F10").Value = 234442875
$ws.Range("G10").Value = 6.35501

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "DOGE"
$ws.Range("C11").Value = "Dogecoin"
$ws.Range("D11").Value = 0.070879
$ws.Range("E11").Value = 10044238864
$ws.Range("F11").Value = 362071857
$ws.Range("G11").Value = 3.24175

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "TON"
$ws.Range("C12").Value = "Toncoin"
$ws.Range("D12").Value = 2.28
$ws.Range("E12").Value = 9064058075
$ws.Range("F12").Value = 12010261
$ws.Range("G12").Value = 1.99037

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "TRX"
$ws.Range("C13").Value = "TRON"
$ws.Range("D13").Value = 0.098971
$ws.Range("E13").Value = 8769367934
$ws.Range("F13").Value = 209642443
$ws.Range("G13").Value = 1.23022

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "LINK"
$ws.Range("C14").Value = "Chainlink"
$ws.Range("D14").Value = 12.15
$ws.Range("E14").Value = 6766887160
$ws.Range("F14").Value = 671965452
$ws.Range("G14").Value = 7.39234

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "MATIC"
$ws.Range("C15").Value = "Polygon"
$ws.Range("D15").Value = 0.6961270000000001
$ws.Range("E15").Value = 6373371461
$ws.Range("F15").Value = 216107257
$ws.Range("G15").Value = 3.80566

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "DOT"
$ws.Range("C16").Value = "Polkadot"
$ws.Range("D16").Value = 4.83
$ws.Range("E16").Value = 6209660529
$ws.Range("F16").Value = 190226600
$ws.Range("G16").Value = 4.04255

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "WBTC"
$ws.Range("C17").Value = "Wrapped Bitcoin"
$ws.Range("D17").Value = 35099
$ws.Range("E17").Value = 5739244211
$ws.Range("F17").Value = 127256614
$ws.Range("G17").Value = 1.01391

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "LTC"
$ws.Range("C18").Value = "Litecoin"
$ws.Range("D18").Value = 71.66
$ws.Range("E18").Value = 5291206154
$ws.Range("F18").Value = 495598951
$ws.Range("G18").Value = 2.90199

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "SHIB"
$ws.Range("C19").Value = "Shiba Inu"
$ws.Range("D19").Value = 0.00000817
$ws.Range("E19").Value = 4800231104
$ws.Range("F19").Value = 167032983
$ws.Range("G19").Value = 3.75247

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "BCH"
$ws.Range("C20").Value = "Bitcoin Cash"
$ws.Range("D20").Value = 242.13
$ws.Range("E20").Value = 4733473724
$ws.Range("F20").Value = 94741860
$ws.Range("G20").Value = 2.09059

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "AVAX"
$ws.Range("C21").Value = "Avalanche"
$ws.Range("D21").Value = 12.57
$ws.Range("E21").Value = 4466763899
$ws.Range("F21").Value = 189783815
$ws.Range("G21").Value = 4.57189

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "DAI"
$ws.Range("C22").Value = "Dai"
$ws.Range("D22").Value = 0.999287
$ws.Range("E22").Value = 3722563356
$ws.Range("F22").Value = 117423740
$ws.Range("G22").Value = 0.04194

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "LEO Token"
$ws.Range("D23").Value = 3.94
$ws.Range("E23").Value = 3655689501
$ws.Range("F23").Value = 534110
$ws.Range("G23").Value = -0.64447

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "UNI"
$ws.Range("C24").Value = "Uniswap"
$ws.Range("D24").Value = 4.79
$ws.Range("E24").Value = 3598636704
$ws.Range("F24").Value = 320213940
$ws.Range("G24").Value = 1.26873

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "XLM"
$ws.Range("C25").Value = "Stellar"
$ws.Range("D25").Value = 0.126972
$ws.Range("E25").Value = 3533750066
$ws.Range("F25").Value = 84418258
$ws.Range("G25").Value = 3.18377

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "TUSD"
$ws.Range("C26").Value = "TrueUSD"
$ws.Range("D26").Value = 0.999024
$ws.Range("E26").Value = 3332262803
$ws.Range("F26").Value = 132869399
$ws.Range("G26").Value = -0.03853

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "XMR"
$ws.Range("C27").Value = "Monero"
$ws.Range("D27").Value = 169.47
$ws.Range("E27").Value = 3073093932
$ws.Range("F27").Value = 53362385
$ws.Range("G27").Value = 0.24299

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "OKB"
$ws.Range("C28").Value = "OKB"
$ws.Range("D28").Value = 48.66
$ws.Range("E28").Value = 2913254912
$ws.Range("F28").Value = 16979495
$ws.Range("G28").Value = 4.78792

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "ETC"
$ws.Range("C29").Value = "Ethereum Classic"
$ws.Range("D29").Value = 18.18
$ws.Range("E29").Value = 2599938731
$ws.Range("F29").Value = 83208055
$ws.Range("G29").Value = 4.76539

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "ATOM"
$ws.Range("C30").Value = "Cosmos Hub"
$ws.Range("D30").Value = 8.390000000000001
$ws.Range("E30").Value = 2454193351
$ws.Range("F30").Value = 201724445
$ws.Range("G30").Value = 6.55821

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "FIL"
$ws.Range("C31").Value = "Filecoin"
$ws.Range("D31").Value = 4.11
$ws.Range("E31").Value = 1906661475
$ws.Range("F31").Value = 115763225
$ws.Range("G31").Value = 3.90937

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "BUSD"
$ws.Range("C32").Value = "BUSD"
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 1897531359
$ws.Range("F32").Value = 3235136233
$ws.Range("G32").Value = -0.0323

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "HBAR"
$ws.Range("C33").Value = "Hedera"
$ws.Range("D33").Value = 0.056376
$ws.Range("E33").Value = 1890777818
$ws.Range("F33").Value = 25487278
$ws.Range("G33").Value = 2.65749

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "ICP"
$ws.Range("C34").Value = "Internet Computer"
$ws.Range("D34").Value = 4.15
$ws.Range("E34").Value = 1861459510
$ws.Range("F34").Value = 30287938
$ws.Range("G34").Value = 4.12787

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "CRO"
$ws.Range("C35").Value = "Cronos"
$ws.Range("D35").Value = 0.06884999999999999
$ws.Range("E35").Value = 1813222675
$ws.Range("F35").Value = 9452427
$ws.Range("G35").Value = 2.2761

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "LDO"
$ws.Range("C36").Value = "Lido DAO"
$ws.Range("D36").Value = 2.03
$ws.Range("E36").Value = 1803670925
$ws.Range("F36").Value = 48772634
$ws.Range("G36").Value = 1.18192

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "APT"
$ws.Range("C37").Value = "Aptos"
$ws.Range("D37").Value = 7.01
$ws.Range("E37").Value = 1742531649
$ws.Range("F37").Value = 77872833
$ws.Range("G37").Value = 1.4844

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "NEAR"
$ws.Range("C38").Value = "NEAR Protocol"
$ws.Range("D38").Value = 1.63
$ws.Range("E38").Value = 1608683341
$ws.Range("F38").Value = 139110282
$ws.Range("G38").Value = 9.21252

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "VET"
$ws.Range("C39").Value = "VeChain"
$ws.Range("D39").Value = 0.02078443
$ws.Range("E39").Value = 1510140826
$ws.Range("F39").Value = 41336559
$ws.Range("G39").Value = 3.69558

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "QNT"
$ws.Range("C40").Value = "Quant"
$ws.Range("D40").Value = 101.31
$ws.Range("E40").Value = 1467701750
$ws.Range("F40").Value = 19746836
$ws.Range("G40").Value = 1.54054

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "ARB"
$ws.Range("C41").Value = "Arbitrum"
$ws.Range("D41").Value = 1.1
$ws.Range("E41").Value = 1407878860
$ws.Range("F41").Value = 301587224
$ws.Range("G41").Value = 3.11647

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "OP"
$ws.Range("C42").Value = "Optimism"
$ws.Range("D42").Value = 1.54
$ws.Range("E42").Value = 1349299437
$ws.Range("F42").Value = 126616269
$ws.Range("G42").Value = 6.16335

$ws.Range("A43").Value = 43
$ws.Range("B43").Value = "AAVE"
$ws.Range("C43").Value = "Aave"
$ws.Range("D43").Value = 91.31
$ws.Range("E43").Value = 1336074956
$ws.Range("F43").Value = 124207911
$ws.Range("G43").Value = 0.72412

$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "KAS"
$ws.Range("C44").Value = "Kaspa"
$ws.Range("D44").Value = 0.062127
$ws.Range("E44").Value = 1329531275
$ws.Range("F44").Value = 44839987
$ws.Range("G44").Value = 13.29736

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "INJ"
$ws.Range("C45").Value = "Injective"
$ws.Range("D45").Value = 15.65
$ws.Range("E45").Value = 1311769081
$ws.Range("F45").Value = 68968162
$ws.Range("G45").Value = 6.78456

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "MNT"
$ws.Range("C46").Value = "Mantle"
$ws.Range("D46").Value = 0.414792
$ws.Range("E46").Value = 1288338765
$ws.Range("F46").Value = 56708987
$ws.Range("G46").Value = 1.24488

$ws.Range("A47").Value = 47
$ws.Range("B47").Value = "EGLD"
$ws.Range("C47").Value = "MultiversX"
$ws.Range("D47").Value = 46.43
$ws.Range("E47").Value = 1216981858
$ws.Range("F47").Value = 290287143
$ws.Range("G47").Value = 37.59187

$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "MKR"
$ws.Range("C48").Value = "Maker"
$ws.Range("D48").Value = 1339.55
$ws.Range("E48").Value = 1207352218
$ws.Range("F48").Value = 33851538
$ws.Range("G48").Value = 0.133

$ws.Range("A49").Value = 49
$ws.Range("B49").Value = "IMX"
$ws.Range("C49").Value = "ImmutableX"
$ws.Range("D49").Value = 0.951439
$ws.Range("E49").Value = 1190900088
$ws.Range("F49").Value = 770261663
$ws.Range("G49").Value = 24.78934

$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "GRT"
$ws.Range("C50").Value = "The Graph"
$ws.Range("D50").Value = 0.12779
$ws.Range("E50").Value = 1184634620
$ws.Range("F50").Value = 125277039
$ws.Range("G50").Value = 9.78881

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "RETH"
$ws.Range("C51").Value = "Rocket Pool ETH"
$ws.Range("D51").Value = 2072.79
$ws.Range("E51").Value = 1109107345
$ws.Range("F51").Value = 8102440
$ws.Range("G51").Value = 2.49026
